$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 107
$ws.Range("C2").Value = "face/face007.jpg"
$ws.Range("D2").Value = "proben"
$ws.Range("E2").Value = "face"

$ws.Range("B3").Value = 34
$ws.Range("C3").Value = "flower/flower001.jpg"
$ws.Range("D3").Value = "zielen"
$ws.Range("E3").Value = "flower"

$ws.Range("B4").Value = 120
$ws.Range("C4").Value = "flower/flower014.jpg"
$ws.Range("D4").Value = "ändern"
$ws.Range("E4").Value = "flower"

$ws.Range("B5").Value = 67
$ws.Range("C5").Value = "face/face026.jpg"
$ws.Range("D5").Value = "lächeln"
$ws.Range("E5").Value = "face"

$ws.Range("B6").Value = 72
$ws.Range("C6").Value = "flower/flower016.jpg"
$ws.Range("D6").Value = "quellen"
$ws.Range("E6").Value = "flower"

$ws.Range("B7").Value = 121
$ws.Range("C7").Value = "face/face006.jpg"
$ws.Range("D7").Value = "tollen"
$ws.Range("E7").Value = "face"

$ws.Range("B8").Value = 117
$ws.Range("C8").Value = "flower/flower012.jpg"
$ws.Range("D8").Value = "rechnen"
$ws.Range("E8").Value = "flower"

$ws.Range("B9").Value = 126
$ws.Range("C9").Value = "face/face012.jpg"
$ws.Range("D9").Value = "kennen"
$ws.Range("E9").Value = "face"

$ws.Range("B10").Value = 44
$ws.Range("C10").Value = "face/face023.jpg"
$ws.Range("D10").Value = "achten"
$ws.Range("E10").Value = "face"

$ws.Range("B11").Value = 114
$ws.Range("C11").Value = "flower/flower028.jpg"
$ws.Range("D11").Value = "reisen"
$ws.Range("E11").Value = "flower"

$ws.Range("B12").Value = 93
$ws.Range("C12").Value = "flower/flower023.jpg"
$ws.Range("D12").Value = "lassen"
$ws.Range("E12").Value = "flower"

$ws.Range("B13").Value = 97
$ws.Range("C13").Value = "flower/flower008.jpg"
$ws.Range("D13").Value = "öffnen"
$ws.Range("E13").Value = "flower"

$ws.Range("B14").Value = 15
$ws.Range("C14").Value = "face/face011.jpg"
$ws.Range("D14").Value = "frischen"
$ws.Range("E14").Value = "face"

$ws.Range("B15").Value = 89
$ws.Range("C15").Value = "flower/flower004.jpg"
$ws.Range("D15").Value = "bergen"
$ws.Range("E15").Value = "flower"

$ws.Range("B16").Value = 98
$ws.Range("C16").Value = "face/face014.jpg"
$ws.Range("D16").Value = "herrschen"
$ws.Range("E16").Value = "face"

$ws.Range("B17").Value = 55
$ws.Range("C17").Value = "face/face020.jpg"
$ws.Range("D17").Value = "sparen"
$ws.Range("E17").Value = "face"

$ws.Range("B18").Value = 43
$ws.Range("C18").Value = "face/face010.jpg"
$ws.Range("D18").Value = "leeren"
$ws.Range("E18").Value = "face"

$ws.Range("B19").Value = 42
$ws.Range("C19").Value = "flower/flower025.jpg"
$ws.Range("D19").Value = "holen"
$ws.Range("E19").Value = "flower"

$ws.Range("B20").Value = 90
$ws.Range("C20").Value = "face/face005.jpg"
$ws.Range("D20").Value = "fließen"
$ws.Range("E20").Value = "face"

$ws.Range("B21").Value = 54
$ws.Range("C21").Value = "face/face002.jpg"
$ws.Range("D21").Value = "wachsen"
$ws.Range("E21").Value = "face"

$ws.Range("B22").Value = 50
$ws.Range("C22").Value = "flower/flower015.jpg"
$ws.Range("D22").Value = "atmen"
$ws.Range("E22").Value = "flower"

$ws.Range("B23").Value = 91
$ws.Range("C23").Value = "flower/flower030.jpg"
$ws.Range("D23").Value = "deuten"
$ws.Range("E23").Value = "flower"

$ws.Range("B24").Value = 84
$ws.Range("C24").Value = "face/face024.jpg"
$ws.Range("D24").Value = "danken"
$ws.Range("E24").Value = "face"

$ws.Range("B25").Value = 8
$ws.Range("C25").Value = "flower/flower005.jpg"
$ws.Range("D25").Value = "planen"
$ws.Range("E25").Value = "flower"

$ws.Range("B26").Value = 113
$ws.Range("C26").Value = "flower/flower027.jpg"
$ws.Range("D26").Value = "trotzen"
$ws.Range("E26").Value = "flower"

$ws.Range("B27").Value = 22
$ws.Range("C27").Value = "flower/flower019.jpg"
$ws.Range("D27").Value = "stoppen"
$ws.Range("E27").Value = "flower"

$ws.Range("B28").Value = 53
$ws.Range("C28").Value = "face/face027.jpg"
$ws.Range("D28").Value = "wecken"
$ws.Range("E28").Value = "face"

$ws.Range("B29").Value = 124
$ws.Range("C29").Value = "face/face008.jpg"
$ws.Range("D29").Value = "angeln"
$ws.Range("E29").Value = "face"

$ws.Range("B30").Value = 0
$ws.Range("C30").Value = "face/face022.jpg"
$ws.Range("D30").Value = "prüfen"
$ws.Range("E30").Value = "face"

$ws.Range("B31").Value = 32
$ws.Range("C31").Value = "flower/flower006.jpg"
$ws.Range("D31").Value = "bauen"
$ws.Range("E31").Value = "flower"

$ws.Range("B32").Value = 52
$ws.Range("C32").Value = "flower/flower020.jpg"
$ws.Range("D32").Value = "kriegen"
$ws.Range("E32").Value = "flower"

$ws.Range("B33").Value = 109
$ws.Range("C33").Value = "face/face013.jpg"
$ws.Range("D33").Value = "wehen"
$ws.Range("E33").Value = "face"
